$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old data rows (rows 3 through 13), keeping header row 1 and first data row 2
$ws.Range("A3:B13").EntireRow.Delete() | Out-Null

# Replace row 2 contents with the new, more realistic demo data
$ws.Range("B2").Value = "ofs-pop"
$ws.Range("A2").Value = "population"

# Shrink the table / list object to match the new data extent
$ws.ListObjects.Item(1).Resize($ws.Range("A1:B2"))

# Column A is now slightly wider to fit "population"
$ws.Columns.Item(1).ColumnWidth = 9

# Update the active selection to reflect where the user ended up after editing
$ws.Range("B6").Select() | Out-Null
